# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F4").Value = 970
    $ws.Range("F6").Value = 2341
    $ws.Range("F8").Value = 1433
    $ws.Range("F12").Value = 387
}
